$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Sheet: Restricciones_del_lider
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("A2") "1.9 - x"
Set-TextValue $ws2.Range("B2") "-2.9"
Set-TextValue $ws2.Range("D2") "0.83"
Set-TextValue $ws2.Range("A3") "-1.9 + x"
Set-TextValue $ws2.Range("B3") "0.8999999999999999"
Set-TextValue $ws2.Range("D3") "0.08"

# Sheet: Restricciones_del_follower
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("A2") "-0.37951807228915646 + 2.5301204819277103y"
Set-TextValue $ws3.Range("B2") "-0.6204819277108435"
Set-TextValue $ws3.Range("D2") "0.01"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "6.3"
Set-TextValue $ws3.Range("A3") "1.1102230246251565e-16y"
Set-TextValue $ws3.Range("B3") "-1.0"
Set-TextValue $ws3.Range("D3") "0.97"
Set-TextValue $ws3.Range("E3") "3.9000000000000004"
Set-TextValue $ws3.Range("F3") "0"

# Sheet: Punto_modificado
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "1.9"
Set-TextValue $ws4.Range("B2") "0.15"

# Sheet: Vector_bf
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "0.5379487951807227"

# Sheet: Vector_BF
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "2.347"
Set-TextValue $ws6.Range("A3") "4.225849999999999"

# Sheet: Vector_Alpha (this one is a real numeric cell, not text)
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.4899999999999998
